$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr)
$newRows = @(
    @(" Abu Dhabi", " October 30 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Kings XI Punjab", "Rajasthan Royals", "Deepak Hooda ", "1", "1", "0", "0", "100.00"),
    @(" Abu Dhabi", " November 01 2020", "Super Kings won by 9 wickets (with 7 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "Deepak Hooda ", "62", "30", "3", "4", "206.66"),
    @(" Dubai (DSC)", " October 20 2020", "Kings XI won by 5 wickets (with 6 balls remaining)", "Kings XI Punjab", "Delhi Capitals", "Deepak Hooda ", "15", "22", "1", "0", "68.18"),
    @(" Dubai (DSC)", " October 18 2020", "Match tied (Kings XI won the one-over eliminator)", "Kings XI Punjab", "Mumbai Indians", "Deepak Hooda ", "23", "16", "1", "1", "143.75"),
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Kings XI Punjab", "Sunrisers Hyderabad", "Deepak Hooda ", "0", "2", "0", "0", "0.00")
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        # Force text storage so numeric-looking strings (runs, balls, 4s, 6s, sr)
        # stay text rather than being coerced to numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col - 1]
    }
}
